# Insert a new data row for "Haba" (Vega Modelo de Temuco) at row 54,
# pushing the former rows 54-71 down to 55-72, and populate the new
# row with its own values (weekly update to the consolidated price series).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 54..71 down by one to make room for the new record.
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new weekly record.
$ws.Cells.Item(54, 1).Value = 10
$ws.Cells.Item(54, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(54, 3).Value = "La Araucanía"
$ws.Cells.Item(54, 4).Value = 44798
$ws.Cells.Item(54, 5).Value = 9
$ws.Cells.Item(54, 6).Value = 100112026
$ws.Cells.Item(54, 7).Value = "Haba"
$ws.Cells.Item(54, 8).Value = "Sin especificar"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 70
$ws.Cells.Item(54, 11).Value = 13000
$ws.Cells.Item(54, 12).Value = 15000
$ws.Cells.Item(54, 13).Value = 14143
$ws.Cells.Item(54, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(54, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(54, 16).Value = 566
$ws.Cells.Item(54, 17).Value = 25
$ws.Cells.Item(54, 18).Value = "Hortaliza"
